$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the "Good" cell style to row 6 and row 7 (A/B/C), matching rows 1-5,
# and make column B text-formatted like the other "Good" rows.
$ws.Range("A6").Style = "Good"
$ws.Range("B6").Style = "Good"
$ws.Range("B6").NumberFormat = "@"
$ws.Range("C6").Style = "Good"

$ws.Range("A7").Style = "Good"
$ws.Range("B7").Style = "Good"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("C7").Style = "Good"

# Row 8 / Row 9: new master entries. Reset B8/B9 back to the plain "Normal"
# style (they previously only carried the text-number-format style) and add
# the new A/C values for the new masters.
$ws.Range("B8").Style = "Normal"
$ws.Range("A8").Value = "700-799"
$ws.Range("C8").Value = "Andrew"

$ws.Range("B9").Style = "Normal"
$ws.Range("A9").Value = "800-899"
$ws.Range("C9").Value = "Miho"

# Move the active selection from L6 to M6.
[void]$ws.Range("M6").Select()
